# 12/27 - Undercarriage and skirts work
# - Finished rear exhaust and tagging pieces with IDs for the skirt.
# - Publishing 104-111 and reserving 112 and 113 for the display.
# - Fixed up the electronics undercarriage file with parts already designed.
# - Added PSU to the model and designed a power supply mount, as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row at sheet row 62 (pushes the old "90 - Misc - Bracket Cover"
# row, and everything after it, down by one) so the new part #86 (PSU Mount)
# can be slotted in at its sorted position.
$ws.Rows.Item(62).Insert()

# Row 62: part 86 - Electrical - Mount - PSU Mount
$ws.Cells.Item(62,1).Value = 86
$ws.Cells.Item(62,2).Value = "Electrical"
$ws.Cells.Item(62,3).Value = "Mount"
$ws.Cells.Item(62,4).Value = "N"
$ws.Cells.Item(62,5).Value = "PSU Mount"
$ws.Cells.Item(62,6).Value = "ABS"
$ws.Cells.Item(62,7).Value = 2
$ws.Cells.Item(62,9).Value = "86 - Electrical - Mount - PSU Mount.stl"

# Row 76: part 104 - Exterior - Skirt - L Front End
$ws.Cells.Item(76,1).Value = 104
$ws.Cells.Item(76,2).Value = "Exterior"
$ws.Cells.Item(76,3).Value = "Skirt"
$ws.Cells.Item(76,4).Value = "L"
$ws.Cells.Item(76,5).Value = "L Front End"
$ws.Cells.Item(76,6).Value = "ABS"
$ws.Cells.Item(76,7).Value = 1
$ws.Cells.Item(76,9).Value = "104 - Exterior - Skirt - L Front End.stl"

# Row 77: part 105 - Exterior - Skirt - R Front End
$ws.Cells.Item(77,1).Value = 105
$ws.Cells.Item(77,2).Value = "Exterior"
$ws.Cells.Item(77,3).Value = "Skirt"
$ws.Cells.Item(77,4).Value = "R"
$ws.Cells.Item(77,5).Value = "R Front End"
$ws.Cells.Item(77,6).Value = "ABS"
$ws.Cells.Item(77,7).Value = 1
$ws.Cells.Item(77,8).Value = "Power switch"
$ws.Cells.Item(77,9).Value = "105 - Exterior - Skirt - R Front End.stl"

# Row 78: part 106 - Exterior - Skirt - L Rear End
$ws.Cells.Item(78,1).Value = 106
$ws.Cells.Item(78,2).Value = "Exterior"
$ws.Cells.Item(78,3).Value = "Skirt"
$ws.Cells.Item(78,4).Value = "L"
$ws.Cells.Item(78,5).Value = "L Rear End"
$ws.Cells.Item(78,6).Value = "ABS"
$ws.Cells.Item(78,7).Value = 1
$ws.Cells.Item(78,9).Value = "106 - Exterior - Skirt - L Rear End.stl"

# Row 79: part 107 - Exterior - Skirt - R Rear End
$ws.Cells.Item(79,1).Value = 107
$ws.Cells.Item(79,2).Value = "Exterior"
$ws.Cells.Item(79,3).Value = "Skirt"
$ws.Cells.Item(79,4).Value = "R"
$ws.Cells.Item(79,5).Value = "R Rear End"
$ws.Cells.Item(79,6).Value = "ABS"
$ws.Cells.Item(79,7).Value = 1
$ws.Cells.Item(79,8).Value = "IEC (power) plug"
$ws.Cells.Item(79,9).Value = "107 - Exterior - Skirt - R Rear End.stl"

# Row 80: part 108 - Exterior - Skirt - Grille
$ws.Cells.Item(80,1).Value = 108
$ws.Cells.Item(80,2).Value = "Exterior"
$ws.Cells.Item(80,3).Value = "Skirt"
$ws.Cells.Item(80,4).Value = "N"
$ws.Cells.Item(80,5).Value = "Grille"
$ws.Cells.Item(80,6).Value = "ABS"
$ws.Cells.Item(80,7).Value = 4
$ws.Cells.Item(80,9).Value = "108 - Exterior - Skirt - Grille.stl"

# Row 81: part 109 - Exterior - Skirt - 92mm Fan Mount
$ws.Cells.Item(81,1).Value = 109
$ws.Cells.Item(81,2).Value = "Exterior"
$ws.Cells.Item(81,3).Value = "Skirt"
$ws.Cells.Item(81,4).Value = "N"
$ws.Cells.Item(81,5).Value = "92mm Fan Mount"
$ws.Cells.Item(81,6).Value = "ABS"
$ws.Cells.Item(81,7).Value = 3
$ws.Cells.Item(81,9).Value = "109 - Exterior - Skirt - 92mm Fan Mount.stl"

# Row 82: part 110 - Exterior - Skirt - Fan Intake
$ws.Cells.Item(82,1).Value = 110
$ws.Cells.Item(82,2).Value = "Exterior"
$ws.Cells.Item(82,3).Value = "Skirt"
$ws.Cells.Item(82,4).Value = "N"
$ws.Cells.Item(82,5).Value = "Fan Intake"
$ws.Cells.Item(82,6).Value = "ABS"
$ws.Cells.Item(82,7).Value = 1
$ws.Cells.Item(82,8).Value = "Interchangeable with 111"
$ws.Cells.Item(82,9).Value = "110 - Exterior - Skirt - Fan Intake.stl"

# Row 83: part 111 - Exterior - Skirt - Fan Exhaust
$ws.Cells.Item(83,1).Value = 111
$ws.Cells.Item(83,2).Value = "Exterior"
$ws.Cells.Item(83,3).Value = "Skirt"
$ws.Cells.Item(83,4).Value = "N"
$ws.Cells.Item(83,5).Value = "Fan Exhuast"
$ws.Cells.Item(83,6).Value = "ABS"
$ws.Cells.Item(83,7).Value = 2
$ws.Cells.Item(83,8).Value = "Interchangeable with 110"
$ws.Cells.Item(83,9).Value = "111 - Exterior - Skirt - Fan Exhaust.stl"

# Row 84: part 112 - Exterior - Skirt - Screen Mount (reserved)
$ws.Cells.Item(84,1).Value = 112
$ws.Cells.Item(84,2).Value = "Exterior"
$ws.Cells.Item(84,3).Value = "Skirt"
$ws.Cells.Item(84,4).Value = "N"
$ws.Cells.Item(84,5).Value = "Screen Mount"
$ws.Cells.Item(84,6).Value = "ABS"
$ws.Cells.Item(84,7).Value = 1
$ws.Cells.Item(84,9).Value = "112 - Exterior - Skirt - Screen Mount.stl"

# Row 85: part 113 - Exterior - Skirt - Screen Adapter (reserved)
$ws.Cells.Item(85,1).Value = 113
$ws.Cells.Item(85,2).Value = "Exterior"
$ws.Cells.Item(85,3).Value = "Skirt"
$ws.Cells.Item(85,4).Value = "N"
$ws.Cells.Item(85,5).Value = "Screen Adapter"
$ws.Cells.Item(85,6).Value = "ABS"
$ws.Cells.Item(85,7).Value = 1
$ws.Cells.Item(85,9).Value = "113 - Exterior - Skirt - Screen Adapter.stl"

# Grow the table to cover the newly inserted row and the appended rows.
$lo.Resize($ws.Range("A1:I85"))

# Leave the selection where the author's last edit landed.
$ws.Range("F77").Select()
